$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new "Resumen de Reunión" / "Revisión de SQA" comments in column J ---
# Shared strings must be created in this order (J10, J9, J11) so that the
# underlying shared-string table receives new unique entries in the same
# order as the target workbook (J10 -> 147, J9 -> 148, J11 -> 149).
$ws.Range("J10").Value = "Resumen de Reunión 23"
$ws.Range("J9").Value = "Resumen de Reunión 22"
$ws.Range("J11").Value = "Revisión de SQA - Arquitectura del sistema"

# --- Widen column J slightly to better fit the new text ---
$ws.Columns.Item(10).ColumnWidth = 38

# --- Update the view: scroll so column D is visible and select K14 ---
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("K14").Select()
